$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at row 36 (pushes old rows 36+ down by one) ---
$ws.Rows(36).Insert()

# --- 2. Row 36 should inherit the "last row" special bottom-border style that
#        row 35 currently has (row 35 was previously the last data row). ---
$ws.Range("B35:J35").Copy()
$ws.Range("B36:J36").PasteSpecial(-4122)  # xlPasteFormats

# --- 3. Row 35 is no longer the last row, so it should take on the "normal"
#        style used by the other middle rows (copy from row 34). ---
$ws.Range("B34:J34").Copy()
$ws.Range("B35:J35").PasteSpecial(-4122)  # xlPasteFormats

# --- 4. Fill in the new row 36 data (period 2509, same worker as other rows) ---
$ws.Range("B36").Value = "CC"
$ws.Range("C36").Value = "1007275831"
$ws.Range("D36").Value = "EVA YULIANIS GUERRA VARGAS"
$ws.Range("E36").Value = "2509"
$ws.Range("F36").Value = 52000
$ws.Range("G36").Value = 1300000

# --- 5. Column E (Periodo Mora) across the data rows now gets centered text ---
$ws.Range("E16:E36").HorizontalAlignment = -4108  # xlCenter

# --- 6. Update the summary figures ---
$ws.Range("E11").Value = 1046933
$ws.Range("F13").Value = 21

Write-Host "Edits applied."
